# Adds three new observation rows (57, 58, 59) to the "Artfynd" sheet,
# following the same column layout / value types as the existing data rows.
#
# Numeric columns (A, B, E, Q, R, S) are written as plain numbers.
# Boolean columns (AD, AE, AG) are written as booleans.
# All other populated columns hold text in the source data (even ones that
# look numeric, e.g. "Antal" in column I), so NumberFormat is forced to "@"
# (Text) before the value is assigned to stop Excel's automatic
# number/date inference from reinterpreting them.
# Columns left blank in the source row (e.g. N, AT, AY, and L on row 58)
# are simply not written to, matching the "empty" cells in those rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($sheet, $row, $col, $text) {
    $cell = $sheet.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

# --- Row 57: Tjäder (Tetrao urogallus) ---
$r = 57
$ws.Cells.Item($r, 1).Value = 131269236          # A Id
$ws.Cells.Item($r, 2).Value = 57073              # B Taxonsorteringsordning
Set-TextCell $ws $r 4 "LC"                       # D Rödlistade
$ws.Cells.Item($r, 5).Value = 100138             # E TaxonId
Set-TextCell $ws $r 6 "Tjäder"                   # F Artnamn
Set-TextCell $ws $r 7 "Tetrao urogallus"         # G Vetenskapligt namn
Set-TextCell $ws $r 8 "Linnaeus, 1758"           # H Auktor
Set-TextCell $ws $r 9 "1"                        # I Antal
Set-TextCell $ws $r 11 "adult"                   # K Ålder-Stadium
Set-TextCell $ws $r 12 "hona"                    # L Kön
Set-TextCell $ws $r 13 "födosökande"             # M Aktivitet
Set-TextCell $ws $r 16 "Risten-Mulstad, Ög"      # P Lokalnamn
$ws.Cells.Item($r, 17).Value = 560702            # Q Ost
$ws.Cells.Item($r, 18).Value = 6463712           # R Nord
$ws.Cells.Item($r, 19).Value = 10                # S Noggrannhet
Set-TextCell $ws $r 20 "Östergötland"            # T Län
Set-TextCell $ws $r 21 "Åtvidaberg"              # U Kommun
Set-TextCell $ws $r 22 "Östergötland"            # V Provins
Set-TextCell $ws $r 23 "Björsäter"               # W Socken
Set-TextCell $ws $r 25 "2026-02-12"              # Y Startdatum
Set-TextCell $ws $r 27 "2026-02-12"              # AA Slutdatum
$ws.Cells.Item($r, 30).Value = $false            # AD Ej återfunnen
$ws.Cells.Item($r, 31).Value = $false            # AE Osäker artbestämning
$ws.Cells.Item($r, 33).Value = $false            # AG Ospontan
Set-TextCell $ws $r 49 "Steve Daurer"            # AW Rapportör
Set-TextCell $ws $r 50 "Steve Daurer"            # AX Observatörer

# --- Row 58: Talltita (Poecile montanus) ---
$r = 58
$ws.Cells.Item($r, 1).Value = 131269217
$ws.Cells.Item($r, 2).Value = 58043
Set-TextCell $ws $r 4 "NT"
$ws.Cells.Item($r, 5).Value = 103021
Set-TextCell $ws $r 6 "Talltita"
Set-TextCell $ws $r 7 "Poecile montanus"
Set-TextCell $ws $r 8 "(Conrad von Baldenstein, 1827)"
Set-TextCell $ws $r 9 "2"
Set-TextCell $ws $r 11 "adult"
Set-TextCell $ws $r 13 "permanent revir"
Set-TextCell $ws $r 16 "Risten-Mulstad, Ög"
$ws.Cells.Item($r, 17).Value = 560808
$ws.Cells.Item($r, 18).Value = 6463646
$ws.Cells.Item($r, 19).Value = 10
Set-TextCell $ws $r 20 "Östergötland"
Set-TextCell $ws $r 21 "Åtvidaberg"
Set-TextCell $ws $r 22 "Östergötland"
Set-TextCell $ws $r 23 "Björsäter"
Set-TextCell $ws $r 25 "2026-02-12"
Set-TextCell $ws $r 27 "2026-02-12"
Set-TextCell $ws $r 29 "Revirparet"              # AC Publik kommentar
$ws.Cells.Item($r, 30).Value = $false
$ws.Cells.Item($r, 31).Value = $false
$ws.Cells.Item($r, 33).Value = $false
Set-TextCell $ws $r 49 "Steve Daurer"
Set-TextCell $ws $r 50 "Steve Daurer"

# --- Row 59: Spillkråka (Dryocopus martius) ---
$r = 59
$ws.Cells.Item($r, 1).Value = 131269222
$ws.Cells.Item($r, 2).Value = 57881
Set-TextCell $ws $r 4 "NT"
$ws.Cells.Item($r, 5).Value = 100049
Set-TextCell $ws $r 6 "Spillkråka"
Set-TextCell $ws $r 7 "Dryocopus martius"
Set-TextCell $ws $r 8 "(Linnaeus, 1758)"
Set-TextCell $ws $r 9 "1"
Set-TextCell $ws $r 11 "adult"
Set-TextCell $ws $r 12 "hane"
Set-TextCell $ws $r 13 "permanent revir"
Set-TextCell $ws $r 16 "Risten-Mulstad, Ög"
$ws.Cells.Item($r, 17).Value = 560838
$ws.Cells.Item($r, 18).Value = 6463567
$ws.Cells.Item($r, 19).Value = 10
Set-TextCell $ws $r 20 "Östergötland"
Set-TextCell $ws $r 21 "Åtvidaberg"
Set-TextCell $ws $r 22 "Östergötland"
Set-TextCell $ws $r 23 "Björsäter"
Set-TextCell $ws $r 25 "2026-02-12"
Set-TextCell $ws $r 27 "2026-02-12"
$ws.Cells.Item($r, 30).Value = $false
$ws.Cells.Item($r, 31).Value = $false
$ws.Cells.Item($r, 33).Value = $false
Set-TextCell $ws $r 49 "Steve Daurer"
Set-TextCell $ws $r 50 "Steve Daurer"
